# Add a value to cell A1 on the active sheet (Sheet1) — stored as a shared string.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A1").Value = "version1"
